# Update Work Week and Social Spending
# (commit refreshes the Syria GDP-per-Capita dataset with the newer
#  Clio Infra release, and appends six more years of data: 2011-2016)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- 1. Refresh the "GDP per Capita" values (column E) for the rows whose
#        figure changed in the new data release. The sheet stores these as
#        text (same convention as the rest of the workbook), so each cell is
#        marked as Text before the new value is written - mirrors typing a
#        value into a Text-formatted cell in the Excel UI.
$gdpUpdates = [ordered]@{
    2 = "1084"
    52 = "1403"
    95 = "2072"
    132 = "3838"
    133 = "3609"
    134 = "4441"
    135 = "4916"
    136 = "5504"
    137 = "4844"
    138 = "5592"
    139 = "5781"
    140 = "4844"
    141 = "4881"
    142 = "4819"
    143 = "5050"
    144 = "6051"
    145 = "5855"
    146 = "5797"
    147 = "5598"
    148 = "5004"
    149 = "5246"
    150 = "5270"
    151 = "6059"
    152 = "5643"
    153 = "5992"
    154 = "7243"
    155 = "6403"
    156 = "7685"
    157 = "8878"
    158 = "9526"
    159 = "9094"
    160 = "9561"
    161 = "9580"
    162 = "10374"
    163 = "10984"
    164 = "10817"
    165 = "10581"
    166 = "9792"
    167 = "10026"
    168 = "9200"
    169 = "9055"
    170 = "9913"
    171 = "8735"
    172 = "9087"
    173 = "9424.12259076452"
    174 = "10002.275987573"
    175 = "10072.4865364936"
    176 = "9969.37404509101"
    177 = "9856.48311325427"
    178 = "9517.3721795319"
    179 = "8826.77352349446"
    180 = "8742.97356066853"
    181 = "7960.29181725339"
    182 = "7668.07165178925"
    183 = "7509.3042432845"
    184 = "7529.70161488188"
    185 = "6986.52380066229"
    186 = "7051.42847759348"
    187 = "7037.71822979878"
    188 = "6901.77786540411"
    189 = "6779.53709553894"
    190 = "6600.37871990147"
    191 = "6586.51834804181"
    192 = "6520.61227745012"
}

foreach ($row in $gdpUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $gdpUpdates[$row]
}

# --- 2. Append the newly available years (2011-2016) as new rows.
$newYears = [ordered]@{
    193 = @{ Year = 2011; Value = "5979" }
    194 = @{ Year = 2012; Value = "4829" }
    195 = @{ Year = 2013; Value = "4124" }
    196 = @{ Year = 2014; Value = "3536" }
    197 = @{ Year = 2015; Value = "3049" }
    198 = @{ Year = 2016; Value = "3091" }
}

foreach ($row in $newYears.Keys) {
    $entry = $newYears[$row]
    $ws.Cells.Item($row, 1).Value = 760
    $ws.Cells.Item($row, 2).Value = "Syria"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $entry.Year
    $valueCell = $ws.Cells.Item($row, 5)
    $valueCell.NumberFormat = "@"
    $valueCell.Value = $entry.Value
}
